# LoRA FLOPs plan - resource allocation including evaluation steps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2->3: shift the "3(a) Default LoRA run" row down one, drop Steps from 4000 to 1000 ---
$ws.Range("A2:I2").ClearContents()
$ws.Range("A3").Value = "3(a)"
$ws.Range("B3").Value = "Default LoRA run"
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 0.0001
$ws.Range("E3").Value = 512
$ws.Range("F3").Value = 1000
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 16.44
$ws.Range("I3").Value = 65760

# --- Hyperparam sweep block: cap every "Steps" value at 1000 (clear where the value duplicated the row above) ---
$ws.Range("F5").ClearContents()
$ws.Range("F6").Value = 1000
$ws.Range("F7").ClearContents()
$ws.Range("F8").Value = 1000
$ws.Range("F9").Value = 1000
$ws.Range("F10").Value = 1000
$ws.Range("F11").ClearContents()
$ws.Range("F12").Value = 1000
$ws.Range("F13").ClearContents()

# --- Context run note ---
$ws.Range("K14").Value = "* no 512 since already done"

# --- Final model run: Steps 3000 -> 2000 ---
$ws.Range("F18").Value = 2000

# --- Eval rows: Steps no longer applicable, clear them ---
$ws.Range("F21").ClearContents()
$ws.Range("F22").ClearContents()

# --- New summary block: Total row and Total Budget ---
$ws.Range("K23").Value = "Total Budget"
$ws.Range("K23").Font.Bold = $true

$ws.Range("A24").Value = "Total"
$ws.Range("A24").Font.Bold = $true
$ws.Range("F24").Formula = "=SUM(F3:F22)"
$ws.Range("I24").Formula = "=SUM(I3:I22)"
$ws.Range("K24").Value = 1000000

$ws.Range("D26").Value = "max training steps"
$ws.Range("D26").Font.Bold = $true
$ws.Range("E26").Value = 14776
$ws.Range("E26").Font.Bold = $true

# --- Column width tweaks to fit the new/longer content ---
$ws.Columns.Item(4).EntireColumn.AutoFit()
$ws.Columns.Item(6).EntireColumn.AutoFit()
$ws.Columns.Item(11).EntireColumn.AutoFit()

# --- View: zoom + new selection ---
$ws.Application.ActiveWindow.Zoom = 125
$ws.Range("F22").Select()
